# Changes In HTM Parameters.xlsx - add the "Exp 24" experiment row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of experiment data (row 20), following the same pattern as the
# existing rows: Experiment Folder / Local Area Density / Potential Radius /
# Local-Global Inhibition / NumActiveColumnsPerInhArea / Result Image Name.
$ws.Range("A20").Value = "Exp 23"
$ws.Range("B20").Value = 0.2
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "Local"
$ws.Range("E20").Value = -1
$ws.Range("F20").Value = "Exp 24.png"

# Match the formatting (centered, bordered style) used by the rest of the
# data rows by copying it down from the row above.
$ws.Range("A19:E19").Copy() | Out-Null
$ws.Range("A20:E20").PasteSpecial(-4122) | Out-Null

# Reflect the author's final on-screen selection.
$ws.Range("I13").Select() | Out-Null
